$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217, pushing existing rows 217:264 down to 218:265.
$ws.Rows("217:217").Insert()

# Populate the newly inserted row 217 with the new record's data.
$ws.Range("A217").Value = 6
$ws.Range("B217").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C217").Value = "Metropolitana"
$ws.Range("D217").Value = 44508
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112052
$ws.Range("G217").Value = "Albahaca"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 45
$ws.Range("K217").Value = 5000
$ws.Range("L217").Value = 6000
$ws.Range("M217").Value = 5556
$ws.Range("N217").Value = "$/paquete"
$ws.Range("O217").Value = "Región de Arica y Parinacota"
$ws.Range("P217").Value = 5556
$ws.Range("Q217").Value = 1
$ws.Range("R217").Value = "Hortaliza"
